$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H70").Value = 3562.5
$ws.Range("J70").Value = 4000
$ws.Range("L70").Value = 12000
$ws.Range("N70").Value = -12540

$ws.Range("H73").Value = 3562.5
$ws.Range("J73").Value = 4000
$ws.Range("L73").Value = 12000
$ws.Range("N73").Value = -13872

$ws.Range("H135").Value = 1218.6666
$ws.Range("I135").Value = 986.4
$ws.Range("K135").Value = 8877.6
$ws.Range("M135").Value = -6342.6

$ws.Range("H141").Value = 769.2857
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 10992.333
$ws.Range("I36").Value = 10992.333
$ws.Range("K36").Value = 10992.333
$ws.Range("M36").Value = -10646.333

$ws.Range("H122").Value = 9333.333000000001
$ws.Range("I122").Value = 9333.333000000001
$ws.Range("K122").Value = 27999.999
$ws.Range("M122").Value = -25549.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 20010.5
$ws.Range("I33").Value = 20010.5
$ws.Range("K33").Value = 20010.5
$ws.Range("M33").Value = -19674.5

$ws.Range("H86").Value = 2000
$ws.Range("J86").Value = 2000
$ws.Range("L86").Value = 2000
$ws.Range("N86").Value = -4246

$ws.Range("H89").Value = 2000
$ws.Range("J89").Value = 2000
$ws.Range("L89").Value = 10000
$ws.Range("N89").Value = -21232

$ws.Range("H134").Value = 1979.5714
$ws.Range("I134").Value = 2059.5
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 6178.5
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -3643.5
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 1819
$ws.Range("I32").Value = 2064.5715
$ws.Range("J32").Value = 959.5
$ws.Range("K32").Value = 2064.5715
$ws.Range("L32").Value = 959.5
$ws.Range("M32").Value = -1748.5715
$ws.Range("N32").Value = -1591.5

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H86").Value = 29999
$ws.Range("I86").Value = 29999
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 29999
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -28876
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 29999
$ws.Range("I89").Value = 29999
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 149995
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -144379
$ws.Range("N89").ClearContents()

$ws.Range("H134").Value = 1514.1428
$ws.Range("I134").Value = 1514.1428
$ws.Range("K134").Value = 4542.428400000001
$ws.Range("M134").Value = -2007.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1283.1666
$ws.Range("J5").Value = 1999.5
$ws.Range("L5").Value = 5998.5
$ws.Range("N5").Value = -6222.5

$ws.Range("H10").Value = 333.66666
$ws.Range("I10").Value = 9
$ws.Range("J10").Value = 658.3333
$ws.Range("K10").Value = 27
$ws.Range("L10").Value = 1974.9999
$ws.Range("M10").Value = 112
$ws.Range("N10").Value = -2252.9999

$ws.Range("H22").Value = 2778.7222
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 2876.0625
$ws.Range("K22").Value = 6000
$ws.Range("L22").Value = 8628.1875
$ws.Range("M22").Value = -5831
$ws.Range("N22").Value = -8966.1875

$ws.Range("H25").Value = 17882.857
$ws.Range("I25").Value = 90
$ws.Range("J25").Value = 25000
$ws.Range("K25").Value = 270
$ws.Range("L25").Value = 75000
$ws.Range("M25").Value = -101
$ws.Range("N25").Value = -75338

$ws.Range("H27").Value = 2778.7222
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 2876.0625
$ws.Range("K27").Value = 6000
$ws.Range("L27").Value = 8628.1875
$ws.Range("M27").Value = -5898
$ws.Range("N27").Value = -8832.1875

$ws.Range("H30").Value = 17882.857
$ws.Range("I30").Value = 90
$ws.Range("J30").Value = 25000
$ws.Range("K30").Value = 270
$ws.Range("L30").Value = 75000
$ws.Range("M30").Value = -168
$ws.Range("N30").Value = -75204

$ws.Range("H34").Value = 3077.9092
$ws.Range("I34").Value = 986.3333
$ws.Range("J34").Value = 3862.25
$ws.Range("K34").Value = 2958.9999
$ws.Range("L34").Value = 11586.75
$ws.Range("M34").Value = -2874.9999
$ws.Range("N34").Value = -11754.75

$ws.Range("H44").Value = 981.5454999999999
$ws.Range("J44").Value = 1160.7778
$ws.Range("L44").Value = 3482.3334
$ws.Range("N44").Value = -4278.3334

$ws.Range("H46").Value = 6375
$ws.Range("I46").Value = 125
$ws.Range("J46").Value = 9500
$ws.Range("K46").Value = 375
$ws.Range("L46").Value = 28500
$ws.Range("M46").Value = -284
$ws.Range("N46").Value = -28682

$ws.Range("H47").Value = 1002
$ws.Range("I47").Value = 1002
$ws.Range("K47").Value = 3006
$ws.Range("M47").Value = -2575

$ws.Range("H51").Value = 375
$ws.Range("I51").Value = 375
$ws.Range("K51").Value = 1125
$ws.Range("M51").Value = -665

$ws.Range("H59").Value = 1062
$ws.Range("I59").Value = 998.5
$ws.Range("K59").Value = 2995.5
$ws.Range("M59").Value = -2455.5

$ws.Range("H132").Value = 1950
$ws.Range("J132").Value = 1950
$ws.Range("L132").Value = 17550
$ws.Range("N132").Value = -22610

$ws.Range("H135").Value = 1283.1666
$ws.Range("J135").Value = 1999.5
$ws.Range("L135").Value = 17995.5
$ws.Range("N135").Value = -23065.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1055.8334
$ws.Range("I132").Value = 1055.8334
$ws.Range("K132").Value = 3167.5002
$ws.Range("M132").Value = -637.5001999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I40").Value = 3253
$ws.Range("J40").Value = 1252751.2
$ws.Range("K40").Value = 3253
$ws.Range("L40").Value = 1252751.2
$ws.Range("M40").Value = -3117
$ws.Range("N40").Value = -1253023.2

$ws.Range("H46").Value = 1636.091
$ws.Range("I46").Value = 582.8333
$ws.Range("J46").Value = 2900
$ws.Range("K46").Value = 582.8333
$ws.Range("L46").Value = 2900
$ws.Range("M46").Value = -394.8333
$ws.Range("N46").Value = -3276

$ws.Range("H55").Value = 2170.1428
$ws.Range("I55").Value = 800.75
$ws.Range("J55").Value = 3996
$ws.Range("K55").Value = 800.75
$ws.Range("L55").Value = 3996
$ws.Range("M55").Value = -627.75
$ws.Range("N55").Value = -4342

$ws.Range("H68").Value = 1812.375
$ws.Range("I68").Value = 1812.375
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1812.375
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1063.375
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1812.375
$ws.Range("I71").Value = 1812.375
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9061.875
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -5317.875
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 1047.5238
$ws.Range("I82").Value = 795
$ws.Range("J82").Value = 1457.875
$ws.Range("K82").Value = 795
$ws.Range("L82").Value = 1457.875
$ws.Range("M82").Value = -434
$ws.Range("N82").Value = -2179.875

$ws.Range("H85").Value = 1047.5238
$ws.Range("I85").Value = 795
$ws.Range("J85").Value = 1457.875
$ws.Range("K85").Value = 795
$ws.Range("L85").Value = 1457.875
$ws.Range("M85").Value = 453
$ws.Range("N85").Value = -3953.875

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
